$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (pushes the existing data rows down by one).
$ws.Rows.Item(2).Insert()

# Columns whose new value looks like a date string ("YYYY-MM-DD"). Excel's
# auto-detection would otherwise silently convert these into date serial
# numbers; force them to remain plain text by temporarily marking the
# cells as text-formatted, then strip the formatting back off afterwards
# so the cells end up with no explicit style (matching the rest of the
# sheet) while still holding a literal text value.
$dateCells = @("A2", "D2", "E2")
foreach ($addr in $dateCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A2").Value = "2024-04-11"
$ws.Range("B2").Value = "신한제13호스팩"
$ws.Range("C2").Value = "신한"
$ws.Range("D2").Value = "2024-04-15"
$ws.Range("E2").Value = "2024-04-22"
$ws.Range("F2").Value = 6000000
$ws.Range("G2").Value = 3000000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "1337.88 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"

# Drop the explicit text formatting applied above now that the values are
# committed as text, so the cells fall back to the sheet's default style.
$ws.Range("A2:E2").ClearFormats()

# The two oldest records (originally rows 12 & 13: 비엔케이제2호스팩 /
# 유진스팩10호) are no longer part of the feed; after the insert above they
# now live at rows 13 & 14. Remove them so the sheet ends with 12 rows
# (1 header + 11 data rows).
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(13).Delete()
